$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reports")

# Add new row 9: Total label in A9 and total amount in C9
$ws.Range("A9").Value = "Total:"
$ws.Range("C9").Value = 70000
